$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 15048633
$ws.Range("J43").Value = 30096260
$ws.Range("L43").Value = 30096260
$ws.Range("N43").Value = -30096398
$ws.Range("H111").Value = 2482
$ws.Range("I111").Value = 2482
$ws.Range("K111").Value = 7446
$ws.Range("M111").Value = -4379
$ws.Range("H112").Value = 2182.5881
$ws.Range("J112").Value = 2269.625
$ws.Range("L112").Value = 6808.875
$ws.Range("N112").Value = -9024.875
$ws.Range("H127").Value = 1718.1904
$ws.Range("J127").Value = 2537.182
$ws.Range("L127").Value = 7611.545999999999
$ws.Range("N127").Value = -17531.546
$ws.Range("H138").Value = 2833.575
$ws.Range("I138").Value = 2142
$ws.Range("J138").Value = 3095.8965
$ws.Range("K138").Value = 6426
$ws.Range("L138").Value = 9287.6895
$ws.Range("M138").Value = -1286
$ws.Range("N138").Value = -19567.6895
$ws.Range("H141").Value = 1265.7778
$ws.Range("I141").Value = 924
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 2772
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 2408
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1538.931
$ws.Range("I2").Value = 1246.9131
$ws.Range("J2").Value = 2658.3333
$ws.Range("K2").Value = 1246.9131
$ws.Range("L2").Value = 2658.3333
$ws.Range("M2").Value = -1133.9131
$ws.Range("N2").Value = -2884.3333
$ws.Range("H32").Value = 6099.53
$ws.Range("I32").Value = 4838.825
$ws.Range("J32").Value = 11142.35
$ws.Range("K32").Value = 4838.825
$ws.Range("L32").Value = 11142.35
$ws.Range("M32").Value = -4551.825
$ws.Range("N32").Value = -11716.35
$ws.Range("H116").Value = 1538.931
$ws.Range("I116").Value = 1246.9131
$ws.Range("J116").Value = 2658.3333
$ws.Range("K116").Value = 1246.9131
$ws.Range("L116").Value = 2658.3333
$ws.Range("M116").Value = 1047.0869
$ws.Range("N116").Value = -7246.3333
$ws.Range("H132").Value = 2483.4849
$ws.Range("I132").Value = 2298.24
$ws.Range("J132").Value = 3062.375
$ws.Range("K132").Value = 6894.719999999999
$ws.Range("L132").Value = 9187.125
$ws.Range("M132").Value = -4364.719999999999
$ws.Range("N132").Value = -14247.125
$ws.Range("H139").Value = 42478.75
$ws.Range("J139").Value = 42478.75
$ws.Range("L139").Value = 42478.75
$ws.Range("N139").Value = -52758.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1538.931
$ws.Range("I3").Value = 1246.9131
$ws.Range("J3").Value = 2658.3333
$ws.Range("K3").Value = 1246.9131
$ws.Range("L3").Value = 2658.3333
$ws.Range("M3").Value = -1132.9131
$ws.Range("N3").Value = -2886.3333
$ws.Range("H58").Value = 25728.8
$ws.Range("J58").Value = 29661
$ws.Range("L58").Value = 29661
$ws.Range("N58").Value = -30249

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
$ws.Range("H59").Value = 20333.334
$ws.Range("I59").Value = 14500
$ws.Range("J59").Value = 32000
$ws.Range("K59").Value = 14500
$ws.Range("L59").Value = 32000
$ws.Range("M59").Value = -13355
$ws.Range("N59").Value = -34290
$ws.Range("H68").Value = 11375
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 14833.333
$ws.Range("K68").Value = 1000
$ws.Range("L68").Value = 14833.333
$ws.Range("M68").Value = -251
$ws.Range("N68").Value = -16331.333
$ws.Range("H71").Value = 11375
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 14833.333
$ws.Range("K71").Value = 3000
$ws.Range("L71").Value = 44499.999
$ws.Range("M71").Value = 744
$ws.Range("N71").Value = -51987.999
$ws.Range("H74").Value = 33000
$ws.Range("J74").Value = 33000
$ws.Range("L74").Value = 33000
$ws.Range("N74").Value = -34748
$ws.Range("H77").Value = 33000
$ws.Range("J77").Value = 33000
$ws.Range("L77").Value = 99000
$ws.Range("N77").Value = -107736
$ws.Range("H107").Value = 1335
$ws.Range("I107").Value = 697.1429000000001
$ws.Range("K107").Value = 697.1429000000001
$ws.Range("M107").Value = 1222.8571
$ws.Range("H132").Value = 1881.6222
$ws.Range("I132").Value = 1530.3715
$ws.Range("K132").Value = 4591.1145
$ws.Range("M132").Value = -2061.1145
$ws.Range("H141").Value = 358315.34
$ws.Range("J141").Value = 358315.34
$ws.Range("L141").Value = 358315.34
$ws.Range("N141").Value = -368675.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 788.4
$ws.Range("I113").Value = 722
$ws.Range("J113").Value = 832.6667
$ws.Range("K113").Value = 2166
$ws.Range("L113").Value = 2498.0001
$ws.Range("M113").Value = 4
$ws.Range("N113").Value = -6838.0001
$ws.Range("H131").Value = 31297900
$ws.Range("J131").Value = 56678.332
$ws.Range("L131").Value = 170034.996
$ws.Range("N131").Value = -180114.996
$ws.Range("H132").Value = 1251.8422
$ws.Range("I132").Value = 1208
$ws.Range("J132").Value = 1300.5555
$ws.Range("K132").Value = 10872
$ws.Range("L132").Value = 11704.9995
$ws.Range("M132").Value = -8342
$ws.Range("N132").Value = -16764.9995
$ws.Range("H133").Value = 5785.5356
$ws.Range("I133").Value = 3257.5
$ws.Range("J133").Value = 6206.875
$ws.Range("K133").Value = 9772.5
$ws.Range("L133").Value = 18620.625
$ws.Range("M133").Value = -4712.5
$ws.Range("N133").Value = -28740.625
$ws.Range("H139").Value = 3253.5881
$ws.Range("I139").Value = 3265.2
$ws.Range("K139").Value = 9795.599999999999
$ws.Range("M139").Value = -4655.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2014.4565
$ws.Range("I102").Value = 1339.8667
$ws.Range("K102").Value = 1339.8667
$ws.Range("M102").Value = 282.1333
$ws.Range("H113").Value = 1326.4445
$ws.Range("I113").Value = 1362.5714
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1362.5714
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 807.4286
$ws.Range("N113").Value = -5540
$ws.Range("H126").Value = 1848
$ws.Range("I126").Value = 1586.3334
$ws.Range("K126").Value = 4759.0002
$ws.Range("M126").Value = -2289.0002
$ws.Range("H132").Value = 4685.222
$ws.Range("I132").Value = 4871.5
$ws.Range("J132").Value = 3940.111
$ws.Range("K132").Value = 14614.5
$ws.Range("L132").Value = 11820.333
$ws.Range("M132").Value = -12084.5
$ws.Range("N132").Value = -16880.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2592.2666
$ws.Range("I7").Value = 2563.1428
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2563.1428
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2451.1428
$ws.Range("N7").Value = -3224
$ws.Range("H40").Value = 2428.8667
$ws.Range("I40").Value = 2369.0833
$ws.Range("K40").Value = 2369.0833
$ws.Range("M40").Value = -2233.0833
$ws.Range("H61").Value = 1432.25
$ws.Range("I61").Value = 1289.7273
$ws.Range("K61").Value = 1289.7273
$ws.Range("M61").Value = -1087.7273
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H113").Value = 1432.25
$ws.Range("I113").Value = 1289.7273
$ws.Range("K113").Value = 1289.7273
$ws.Range("M113").Value = 880.2727
$ws.Range("H126").Value = 2592.2666
$ws.Range("I126").Value = 2563.1428
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7689.428400000001
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5219.428400000001
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 2616.32
$ws.Range("I132").Value = 2127.8
$ws.Range("J132").Value = 3349.1
$ws.Range("K132").Value = 6383.400000000001
$ws.Range("L132").Value = 10047.3
$ws.Range("M132").Value = -3853.400000000001
$ws.Range("N132").Value = -15107.3
$ws.Range("H135").Value = 34109.875
$ws.Range("J135").Value = 34109.875
$ws.Range("L135").Value = 34109.875
$ws.Range("N135").Value = -44249.875
$ws.Range("H136").Value = 2173.9092
$ws.Range("I136").Value = 1863.5
$ws.Range("J136").Value = 3001.6667
$ws.Range("K136").Value = 5590.5
$ws.Range("L136").Value = 9005.000100000001
$ws.Range("M136").Value = -3040.5
$ws.Range("N136").Value = -14105.0001
$ws.Range("H141").Value = 54371.25
$ws.Range("J141").Value = 53567.145
$ws.Range("L141").Value = 53567.145
$ws.Range("N141").Value = -63927.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H126").Value = 48310130
$ws.Range("I126").Value = 55556516
$ws.Range("J126").Value = 880
$ws.Range("K126").Value = 166669548
$ws.Range("L126").Value = 2640
$ws.Range("M126").Value = -166667078
$ws.Range("N126").Value = -7580
$ws.Range("H132").Value = 4743.968
$ws.Range("I132").Value = 5086.0415
$ws.Range("J132").Value = 3571.1428
$ws.Range("K132").Value = 15258.1245
$ws.Range("L132").Value = 10713.4284
$ws.Range("M132").Value = -12728.1245
$ws.Range("N132").Value = -15773.4284
$ws.Range("H136").Value = 1722.1395
$ws.Range("I136").Value = 674.85
$ws.Range("J136").Value = 2632.8262
$ws.Range("K136").Value = 2024.55
$ws.Range("L136").Value = 7898.4786
$ws.Range("M136").Value = 525.4499999999998
$ws.Range("N136").Value = -12998.4786
$ws.Range("H141").Value = 68300
$ws.Range("J141").Value = 68300
$ws.Range("L141").Value = 68300
$ws.Range("N141").Value = -78660

